# Add data for 2025-10-16
# Updates the 2025 column (L) - and a couple of 2024 column (K) corrections -
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5356
$ws.Range("L3").Value = 5782
$ws.Range("L4").Value = 1413
$ws.Range("L5").Value = 344
$ws.Range("L6").Value = 4796
$ws.Range("L7").Value = 17691

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L6").Value = 66
$ws.Range("L7").Value = 195

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 339
$ws.Range("L3").Value = 407
$ws.Range("L6").Value = 299
$ws.Range("L7").Value = 1171

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L6").Value = 235
$ws.Range("L7").Value = 818

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 203
$ws.Range("L3").Value = 235
$ws.Range("L6").Value = 184
$ws.Range("L7").Value = 676

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 87
$ws.Range("L7").Value = 310

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 148
$ws.Range("L6").Value = 135
$ws.Range("L7").Value = 577
$ws.Range("L8").Value = 1171
$ws.Range("K11").Value = 483
$ws.Range("L11").Value = 292
$ws.Range("L15").Value = 137
$ws.Range("L19").Value = 481
$ws.Range("L20").Value = 437
$ws.Range("L21").Value = 57
$ws.Range("L25").Value = 105
$ws.Range("L29").Value = 995
$ws.Range("L33").Value = 818
$ws.Range("L34").Value = 105
$ws.Range("L37").Value = 676
$ws.Range("L42").Value = 576
$ws.Range("L51").Value = 225
$ws.Range("L52").Value = 356
$ws.Range("L53").Value = 195
$ws.Range("L54").Value = 385
$ws.Range("L55").Value = 184
$ws.Range("K63").Value = 175
$ws.Range("L64").Value = 119
$ws.Range("L67").Value = 605
$ws.Range("L73").Value = 142
$ws.Range("L76").Value = 276
$ws.Range("L84").Value = 173
$ws.Range("L85").Value = 884
$ws.Range("L90").Value = 183
$ws.Range("L94").Value = 218
$ws.Range("L97").Value = 144
$ws.Range("L99").Value = 310
$ws.Range("L101").Value = 17691

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 180
$ws.Range("L6").Value = 137
$ws.Range("L7").Value = 605

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 173

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 95
$ws.Range("L7").Value = 385

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 300
$ws.Range("L3").Value = 381
$ws.Range("L7").Value = 995

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 149
$ws.Range("L7").Value = 481

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 196
$ws.Range("L4").Value = 49
$ws.Range("L6").Value = 160
$ws.Range("L7").Value = 576

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 143
$ws.Range("L7").Value = 437

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 197
$ws.Range("L6").Value = 137
$ws.Range("L7").Value = 577

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 111
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 483
$ws.Range("L7").Value = 292

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L3").Value = 44
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 142

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 148

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 71
$ws.Range("L4").Value = 34
$ws.Range("L7").Value = 225

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 362
$ws.Range("L6").Value = 184
$ws.Range("L7").Value = 884

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L4").Value = 5
$ws.Range("L6").Value = 31

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 117
$ws.Range("L7").Value = 356
